# Fixed issue with driver not setting to null
# - Change the "requireDriverNull"/flag column (C) from "no" to "yes" for a
#   handful of rows on the Sheet1 worksheet.
# - Make Sheet1 the active tab (was InjectSpecificUser) and leave the
#   selection on C17 (was B2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("C2").Value = "yes"
$ws1.Range("C3").Value = "yes"
$ws1.Range("C4").Value = "yes"
$ws1.Range("C14").Value = "yes"
$ws1.Range("C15").Value = "yes"
$ws1.Range("C16").Value = "yes"

$ws1.Activate() | Out-Null
$ws1.Range("C17").Select() | Out-Null
